$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0049205938975016199
$ws.Range("C3").Value = 0.0050545771916707298
$ws.Range("C4").Value = 0.0052857319513956698

$ws.Range("D7").Select()
